$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Rows.Item(13).Insert()
$ws.Range("B12").Value = "the_order"
$ws.Range("B13").Value = "reserved"
$ws.Range("C13").Value = "预留字段"
$ws.Range("H53").Select()
